$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for q0328135_previewuser (old row 5); remaining rows shift up
$ws.Rows.Item(5).Delete() | Out-Null

# Row 2: c1206235
$ws.Cells.Item(2,1).Value2 = "c1206235"
$ws.Cells.Item(2,2).Value2 = "Student ID:  c1206235 `n`n Use the data in https://feb.kuleuven.be/public/u0118298/TASK1/ddKIcCYN/1.dataKIcCYN.txt `n The questions for this task are listed below. `n`n`n Q1: Perform a regression analysis with dependent variable Y2 and independent variables X1 and X3 (without the interaction term). Give the estimate for the parameter corresponding to variable X3. `n`n Q2: Perform a regression analysis with dependent variable Y2 and independent variables X2 and X5. Take the interaction into account. Give the p-value of the parameter corresponding to the interaction term. `n`n Q3: Perform a regression analysis with dependent variable Y3 and independent variables X1, X2 and X3 (without the interaction terms). Give the proportion of explained variance. `n`n`n Don't forget to round decimals to three digits."

# Row 3: c1242115
$ws.Cells.Item(3,1).Value2 = "c1242115"
$ws.Cells.Item(3,2).Value2 = "Cursist ID:  c1242115 `n`n Gebruik de data in https://feb.kuleuven.be/public/u0118298/TASK1/ddEqMRhl/1.dataEqMRhl.txt `n De vragen voor deze taak staan hieronder vermeld. `n`n`n V1: Voer een regressie analyse uit met als afhankelijke variabele Y1 en als onafhankelijke variabelen X1 en X2 (zonder interactieterm). Geef de parameterschatting die hoort bij variabele X2. `n`n V2: Voer een regressie analyse uit met als  afhankelijke variabele Y2, en onafhankelijke variabelen X2 en X5. Neem ook de interactieterm op. Geef de p-waarde die bij de parameter van de interactieterm hoort. `n`n V3: Voer een regressie analyse uit met als afhankelijke variabele Y2, en onafhankelijke variabelen X3, X4 en X5 (zonder interactietermen). Geef de proportie verklaarde variantie. `n`n`n Vergeet kommagetallen niet af te ronden op 3 decimalen."

# Row 4: c1243957
$ws.Cells.Item(4,1).Value2 = "c1243957"
$ws.Cells.Item(4,2).Value2 = "Cursist ID:  c1243957 `n`n Gebruik de data in https://feb.kuleuven.be/public/u0118298/TASK1/ddOVbcNR/1.dataOVbcNR.txt `n De vragen voor deze taak staan hieronder vermeld. `n`n`n V1: Voer een regressie analyse uit met als afhankelijke variabele Y2 en als onafhankelijke variabelen X1 en X3 (zonder interactieterm). Geef de parameterschatting die hoort bij variabele X3. `n`n V2: Voer een regressie analyse uit met als afhankelijke variabele Y1, en onafhankelijke variabelen X2 en X4. Neem ook de interactieterm op. Geef de p-waarde die bij de parameter van de interactieterm hoort. `n`n V3: Voer een regressie analyse uit met als afhankelijke variabele Y1, en onafhankelijke variabelen X2, X3 en X4 (zonder interactietermen). Geef de proportie verklaarde variantie. `n`n`n Vergeet kommagetallen niet af te ronden op 3 decimalen."

# Row 5: q0762379
$ws.Cells.Item(5,1).Value2 = "q0762379"
$ws.Cells.Item(5,2).Value2 = "Student ID:  q0762379 `n`n Use the data in https://feb.kuleuven.be/public/u0118298/TASK1/ddtJEMYS/1.datatJEMYS.txt `n The questions for this task are listed below. `n`n`n Q1: Perform a regression analysis with dependent variable Y1 and independent variables X1 and X2 (without the interaction term). Give the estimate for the parameter corresponding to variable X2. `n`n Q2: Perform a regression analysis with dependent variable Y3 and independent variables X3 and X4. Take the interaction into account. Give the p-value of the parameter corresponding to the interaction term. `n`n Q3: Perform a regression analysis with dependent variable Y3 and independent variables X1, X2 and X3 (without the interaction terms). Give the proportion of explained variance. `n`n`n Don't forget to round decimals to three digits."

# Row 6: q1371623
$ws.Cells.Item(6,1).Value2 = "q1371623"
$ws.Cells.Item(6,2).Value2 = "Student ID:  q1371623 `n`n Use the data in https://feb.kuleuven.be/public/u0118298/TASK1/ddtMfwWh/1.datatMfwWh.txt `n The questions for this task are listed below. `n`n`n Q1: Perform a regression analysis with dependent variable Y1 and independent variables X1 and X2 (without the interaction term). Give the estimate for the parameter corresponding to variable X2. `n`n Q2: Perform a regression analysis with dependent variable Y2 and independent variables X2 and X5. Take the interaction into account. Give the p-value of the parameter corresponding to the interaction term. `n`n Q3: Perform a regression analysis with dependent variable Y3 and independent variables X1, X2 and X3 (without the interaction terms). Give the proportion of explained variance. `n`n`n Don't forget to round decimals to three digits."

# Row 7: q1411379
$ws.Cells.Item(7,1).Value2 = "q1411379"
$ws.Cells.Item(7,2).Value2 = "Student ID:  q1411379 `n`n Use the data in https://feb.kuleuven.be/public/u0118298/TASK1/ddnfKWJY/1.datanfKWJY.txt `n The questions for this task are listed below. `n`n`n Q1: Perform a regression analysis with dependent variable Y1 and independent variables X1 and X2 (without the interaction term). Give the estimate for the parameter corresponding to variable X2. `n`n Q2: Perform a regression analysis with dependent variable Y2 and independent variables X2 and X5. Take the interaction into account. Give the p-value of the parameter corresponding to the interaction term. `n`n Q3: Perform a regression analysis with dependent variable Y3 and independent variables X1, X2 and X3 (without the interaction terms). Give the proportion of explained variance. `n`n`n Don't forget to round decimals to three digits."

# Restore default (non-custom) row heights, matching the original sheet's formatting
$ws.Range("A2:B7").Rows.AutoFit() | Out-Null

